$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at L (12) for "Engineer Name" -----------------
# This shifts the existing L:R columns (Model Number .. Remarks) one to
# the right, becoming M:S, and carries over their cell styles/formats.
$ws.Range("L1").EntireColumn.Insert()

# New column header / placeholder text
$ws.Range("L1").Value = "Engineer Name"
$ws.Range("L2").Value = "{bookings:eng_name}"

# --- Column widths -------------------------------------------------------
# K (Product) gains an explicit width, L (new Engineer Name) gets its own
# width, and M (Model Number, shifted from old L) also gets a new width.
# (Input values below are chosen so that, after this runtime's internal
# character-width/pixel rounding, the stored width lands as close as
# possible to the target widths of 19.33203125 / 37.109375 / 24.44140625.)
$ws.Columns("K").ColumnWidth = 18.5
$ws.Columns("L").ColumnWidth = 36.33333333333333
$ws.Columns("M").ColumnWidth = 23.666666666666664

# --- View state ------------------------------------------------------------
# Scroll the sheet so column H is the left-most visible column, and select
# N14 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("N14").Select()
